$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 7 (The Enchantress quest) ---
# Rewrite the instructions text (D7): drop the "Note" callout, fix typo
# "ton f time" -> "ton of time", add clause about healing/damage spells,
# bold the creature name and split the trailing sentence into its own paragraph.
$ws.Range("D7").Value = '<p>Lets start learning about spells.</p><p>You have just spent a ton of time crafting weapons and Armour, we will repeat this to craft spells. Spells are useful for caster classes, as well as any one who wants to use them, be they healing or damage based.</p><p>Damage spells and Staves (Two Handed weapons) can raise the characters intelligence, both of which can be bought from the shop or crafted.</p><p>Healing spells are great for characters who want to do Cast and attack or Attack and Cast, like Prophets.</p><p>To get the quest item required, you will need to kill: <strong>Umbering Spirit Lord</strong> on Surface. This creature is further down the list and may require you to upgrade your gear through the shop before being able to take him down. </p><p>This creature has a 15% chance to drop the item, so exploration might be a good choice here.</p>'

# Reduce the level/skill requirements for this quest.
$ws.Range("E7").Value = 90
$ws.Range("G7").Value = 5

# --- Row 8 (Go To Labyrinth quest) ---
# Lower the required level.
$ws.Range("E8").Value = 120

# Split out the battle reward handling: the secondary skill requirement
# (H8/I8: "Quick Feet" / 10) is replaced by a faction requirement instead
# (Q8 required_faction_id / R8 required_faction_level) plus the
# required_game_map_id (S8), all pointing at the renamed "Labyrinth" entry.
$ws.Range("Q8").Value = "Labyrinth"
$ws.Range("R8").Value = 1
$ws.Range("S8").Value = "Labyrinth"
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
